$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - rows 5-13, column F ("想去人数" / want-to-go count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 291
$ws1.Range("F6").Value  = 397
$ws1.Range("F7").Value  = 247
$ws1.Range("F8").Value  = 2327
$ws1.Range("F9").Value  = 387
$ws1.Range("F10").Value = 5811
$ws1.Range("F11").Value = 144
$ws1.Range("F12").Value = 379
$ws1.Range("F13").Value = 3

# Sheet "全部类型" (All types) - same underlying events, offset by one row
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 291
$ws4.Range("F7").Value  = 397
$ws4.Range("F8").Value  = 247
$ws4.Range("F11").Value = 2327
$ws4.Range("F12").Value = 387
$ws4.Range("F13").Value = 5811
$ws4.Range("F14").Value = 144
$ws4.Range("F15").Value = 379
$ws4.Range("F17").Value = 3
